$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text values are purely numeric need to be forced to Text
# format first, otherwise Excel auto-converts them to numbers.
$textRefs = @("S6", "S19", "S21", "S34", "S35", "S37", "S39", "S53")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("T2").Value = 274
$ws.Range("U2").Value = 446
$ws.Range("T3").Value = 136
$ws.Range("U3").Value = 224
$ws.Range("T4").Value = 325
$ws.Range("U4").Value = 95
$ws.Range("T5").Value = 324
$ws.Range("U5").Value = 36
$ws.Range("S6").Value = "23069229"
$ws.Range("T6").Value = 685
$ws.Range("U6").Value = 35
$ws.Range("T7").Value = 30
$ws.Range("U7").Value = 330
$ws.Range("T8").Value = 289
$ws.Range("U8").Value = 127
$ws.Range("T9").Value = 279
$ws.Range("U9").Value = 81
$ws.Range("T10").Value = 100
$ws.Range("U10").Value = 316
$ws.Range("Q11").Value = "PREPARAZIONE MACCHINE CONTEMP"
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 305
$ws.Range("T13").Value = 64
$ws.Range("U13").Value = 296
$ws.Range("T14").Value = 49
$ws.Range("U14").Value = 151
$ws.Range("T15").Value = 143
$ws.Range("U15").Value = 217
$ws.Range("T17").Value = 337
$ws.Range("U17").Value = 303
$ws.Range("T18").Value = 296
$ws.Range("U18").Value = 376
$ws.Range("Q19").Value = "LAVORAZIONE"
$ws.Range("R19").Value = "LU82M  CA3"
$ws.Range("S19").Value = "23028656"
$ws.Range("T19").Value = 64
$ws.Range("U19").Value = 0
$ws.Range("T20").Value = 98
$ws.Range("U20").Value = 622
$ws.Range("R21").Value = "PLS0000024"
$ws.Range("S21").Value = "23028583"
$ws.Range("T21").Value = 224
$ws.Range("U21").Value = 0
$ws.Range("V21").Value = "LG3D 0400"
$ws.Range("W21").Value = 23019359
$ws.Range("T22").Value = 329
$ws.Range("U22").Value = 31
$ws.Range("T23").Value = 105
$ws.Range("U23").Value = 543
$ws.Range("T24").Value = 298
$ws.Range("U24").Value = 62
$ws.Range("T25").Value = 248
$ws.Range("U25").Value = 112
$ws.Range("T26").Value = 190
$ws.Range("U26").Value = 170
$ws.Range("T27").Value = 250
$ws.Range("U27").Value = 110
$ws.Range("T28").Value = 53
$ws.Range("U28").Value = 331
$ws.Range("T29").Value = 128
$ws.Range("U29").Value = 232
$ws.Range("Q30").Value = "LAVORAZIONE"
$ws.Range("T30").Value = 9
$ws.Range("U30").Value = 3
$ws.Range("T32").Value = 324
$ws.Range("U32").Value = 72
$ws.Range("T33").Value = 295
$ws.Range("U33").Value = 101
$ws.Range("Q34").Value = "ATTREZZAGGIO"
$ws.Range("R34").Value = "LSB38004X"
$ws.Range("S34").Value = "23040825"
$ws.Range("T34").Value = 49
$ws.Range("U34").Value = 1
$ws.Range("Q35").Value = "ATTREZZAGGIO"
$ws.Range("R35").Value = "LU74MR FA3"
$ws.Range("S35").Value = "23048860"
$ws.Range("T35").Value = 120
$ws.Range("U35").Value = 0
$ws.Range("W35").Value = 23022816
$ws.Range("T36").Value = 333
$ws.Range("U36").Value = 51
$ws.Range("Q37").Value = "LAVORAZIONE"
$ws.Range("S37").Value = "23066097"
$ws.Range("T37").Value = 350
$ws.Range("U37").Value = 10
$ws.Range("V37").Value = "LENDM1244N-1"
$ws.Range("W37").Value = 23066086
$ws.Range("Q38").Value = "LAVORAZIONE"
$ws.Range("T38").Value = 172
$ws.Range("U38").Value = 8
$ws.Range("Q39").Value = "ATTREZZAGGIO"
$ws.Range("R39").Value = "LU89M  IA3"
$ws.Range("S39").Value = "23059608"
$ws.Range("T39").Value = 52
$ws.Range("U39").Value = 2
$ws.Range("T40").Value = 253
$ws.Range("U40").Value = 35
$ws.Range("Q41").Value = "LAVORAZIONE"
$ws.Range("T41").Value = 330
$ws.Range("U41").Value = 30
$ws.Range("T42").Value = 195
$ws.Range("U42").Value = 165
$ws.Range("T43").Value = 214
$ws.Range("U43").Value = 146
$ws.Range("T44").Value = 39
$ws.Range("U44").Value = 345
$ws.Range("Q45").Value = "LAVORAZIONE"
$ws.Range("T45").Value = 200
$ws.Range("U45").Value = 160
$ws.Range("T46").Value = 63
$ws.Range("U46").Value = 321
$ws.Range("T47").Value = 213
$ws.Range("U47").Value = 147
$ws.Range("T48").Value = 105
$ws.Range("U48").Value = 255
$ws.Range("Q49").Value = "FERMO GENERICO"
$ws.Range("T49").Value = 97
$ws.Range("U49").Value = 263
$ws.Range("T50").Value = 263
$ws.Range("U50").Value = 177
$ws.Range("T51").Value = 33
$ws.Range("U51").Value = 63
$ws.Range("T52").Value = 0
$ws.Range("U52").Value = 33
$ws.Range("S53").Value = "23023197"
$ws.Range("T53").Value = 709
$ws.Range("U53").Value = 11
$ws.Range("W53").Value = 22188204
$ws.Range("T54").Value = 309
$ws.Range("U54").Value = 41
$ws.Range("T55").Value = 405
$ws.Range("U55").Value = 235
$ws.Range("T56").Value = 324
$ws.Range("U56").Value = 36
$ws.Range("T58").Value = 245
$ws.Range("U58").Value = 139
$ws.Range("T60").Value = 181
$ws.Range("U60").Value = 203
$ws.Range("T61").Value = 317
$ws.Range("U61").Value = 67
$ws.Range("T62").Value = 51
$ws.Range("U62").Value = 333
$ws.Range("T63").Value = 340
$ws.Range("U63").Value = 44
$ws.Range("Q64").Value = "LAVORAZIONE"
$ws.Range("T64").Value = 88
$ws.Range("U64").Value = 8
$ws.Range("T65").Value = 245
$ws.Range("U65").Value = 55

Write-Host "Applied 149 cell updates"
